$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 47132.668
$ws.Range("J117").Value = 47132.668
$ws.Range("L117").Value = 47132.668
$ws.Range("N117").Value = -56310.668

$ws.Range("H136").Value = 28945.475
$ws.Range("J136").Value = 28945.475
$ws.Range("L136").Value = 28945.475
$ws.Range("N136").Value = -39145.475

$ws.Range("H137").Value = 3689.309
$ws.Range("I137").Value = 1074.8334
$ws.Range("J137").Value = 4009.449
$ws.Range("K137").Value = 3224.5002
$ws.Range("L137").Value = 12028.347
$ws.Range("M137").Value = -674.5001999999999
$ws.Range("N137").Value = -17128.347

$ws.Range("H139").Value = 33704.168
$ws.Range("J139").Value = 33704.168
$ws.Range("L139").Value = 33704.168
$ws.Range("N139").Value = -43984.168

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1950
$ws.Range("I26").Value = 1950
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1950
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -1620
$ws.Range("N26").ClearContents()

$ws.Range("H32").Value = 20163.139
$ws.Range("I32").Value = 19502.195
$ws.Range("J32").Value = 31002.6
$ws.Range("K32").Value = 19502.195
$ws.Range("L32").Value = 31002.6
$ws.Range("M32").Value = -19215.195
$ws.Range("N32").Value = -31576.6

$ws.Range("H61").Value = 2903.24
$ws.Range("I61").Value = 1547.7037
$ws.Range("K61").Value = 1547.7037
$ws.Range("M61").Value = -1335.7037

$ws.Range("H97").Value = 1385.1666
$ws.Range("I97").Value = 766.6667
$ws.Range("J97").Value = 2003.6666
$ws.Range("K97").Value = 766.6667
$ws.Range("L97").Value = 2003.6666
$ws.Range("M97").Value = -270.6667
$ws.Range("N97").Value = -2995.6666

$ws.Range("H118").Value = 49997.332
$ws.Range("J118").Value = 49997.332
$ws.Range("L118").Value = 49997.332
$ws.Range("N118").Value = -53311.332

$ws.Range("H120").Value = 38642.4
$ws.Range("J120").Value = 38642.4
$ws.Range("L120").Value = 38642.4
$ws.Range("N120").Value = -48318.4

$ws.Range("H136").Value = 2903.24
$ws.Range("I136").Value = 1547.7037
$ws.Range("K136").Value = 4643.1111
$ws.Range("M136").Value = -2093.1111

$ws.Range("H139").Value = 40075.43
$ws.Range("J139").Value = 40075.43
$ws.Range("L139").Value = 40075.43
$ws.Range("N139").Value = -50355.43

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 44726.668
$ws.Range("J59").Value = 44726.668
$ws.Range("L59").Value = 44726.668
$ws.Range("N59").Value = -46420.668

$ws.Range("H81").Value = 8933.333000000001
$ws.Range("J81").Value = 8933.333000000001
$ws.Range("L81").Value = 8933.333000000001
$ws.Range("N81").Value = -11055.333

$ws.Range("H84").Value = 8933.333000000001
$ws.Range("J84").Value = 8933.333000000001
$ws.Range("L84").Value = 26799.999
$ws.Range("N84").Value = -37407.999

$ws.Range("H94").Value = 656.44446
$ws.Range("I94").Value = 586.1177
$ws.Range("J94").Value = 776
$ws.Range("K94").Value = 586.1177
$ws.Range("L94").Value = 776
$ws.Range("M94").Value = -135.1177
$ws.Range("N94").Value = -1678

$ws.Range("H134").Value = 4023.8877
$ws.Range("I134").Value = 1810.88
$ws.Range("K134").Value = 5432.64
$ws.Range("M134").Value = -2897.64

$ws.Range("H137").Value = 55378
$ws.Range("J137").Value = 55378
$ws.Range("L137").Value = 55378
$ws.Range("N137").Value = -65578

$ws.Range("H138").Value = 20523.809
$ws.Range("J138").Value = 20523.809
$ws.Range("L138").Value = 20523.809
$ws.Range("N138").Value = -30803.809

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 48664
$ws.Range("J116").Value = 48664
$ws.Range("L116").Value = 48664
$ws.Range("N116").Value = -57842

$ws.Range("H133").Value = 20561.857
$ws.Range("J133").Value = 20561.857
$ws.Range("L133").Value = 20561.857
$ws.Range("N133").Value = -25621.857

$ws.Range("H134").Value = 1371
$ws.Range("I134").Value = 937.09375
$ws.Range("K134").Value = 2811.28125
$ws.Range("M134").Value = -276.28125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1624.5
$ws.Range("J34").Value = 1719.4667
$ws.Range("L34").Value = 5158.4001
$ws.Range("N34").Value = -5326.4001

$ws.Range("H39").Value = 2157.1428
$ws.Range("J39").Value = 2350
$ws.Range("L39").Value = 7050
$ws.Range("N39").Value = -7638

$ws.Range("H92").Value = 1090.36
$ws.Range("I92").Value = 1046.1666
$ws.Range("J92").Value = 1131.1538
$ws.Range("K92").Value = 3138.4998
$ws.Range("L92").Value = 3393.4614
$ws.Range("M92").Value = -1890.4998
$ws.Range("N92").Value = -5889.4614

$ws.Range("H131").Value = 85272.36
$ws.Range("I131").Value = 8838.25
$ws.Range("J131").Value = 155826.92
$ws.Range("K131").Value = 26514.75
$ws.Range("L131").Value = 467480.76
$ws.Range("M131").Value = -21474.75
$ws.Range("N131").Value = -477560.76

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4988.967
$ws.Range("I97").Value = 1635.5
$ws.Range("J97").Value = 10019.167
$ws.Range("K97").Value = 1635.5
$ws.Range("L97").Value = 10019.167
$ws.Range("M97").Value = -1139.5
$ws.Range("N97").Value = -11011.167

$ws.Range("H113").Value = 1826.0741
$ws.Range("I113").Value = 1832.7333
$ws.Range("J113").Value = 1817.75
$ws.Range("K113").Value = 1832.7333
$ws.Range("L113").Value = 1817.75
$ws.Range("M113").Value = 337.2666999999999
$ws.Range("N113").Value = -6157.75

$ws.Range("H134").Value = 23900
$ws.Range("J134").Value = 23900
$ws.Range("L134").Value = 71700
$ws.Range("N134").Value = -76770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 50672
$ws.Range("J116").Value = 50672
$ws.Range("L116").Value = 50672
$ws.Range("N116").Value = -59850

$ws.Range("H121").Value = 35275.332
$ws.Range("J121").Value = 35275.332
$ws.Range("L121").Value = 35275.332
$ws.Range("N121").Value = -38769.332

$ws.Range("H138").Value = 39344.668
$ws.Range("J138").Value = 39344.668
$ws.Range("L138").Value = 39344.668
$ws.Range("N138").Value = -49624.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 33290.332
$ws.Range("J16").Value = 33290.332
$ws.Range("L16").Value = 33290.332
$ws.Range("N16").Value = -33874.332

$ws.Range("H132").Value = 1706.325
$ws.Range("I132").Value = 1207.9524
$ws.Range("J132").Value = 2257.158
$ws.Range("K132").Value = 3623.857199999999
$ws.Range("L132").Value = 6771.474
$ws.Range("M132").Value = -1093.857199999999
$ws.Range("N132").Value = -11831.474
